$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 215
$ws1.Range("F4").Value = 11
$ws1.Range("F5").Value = 6606
$ws1.Range("F9").Value = 5989
$ws1.Range("F16").Value = 92
$ws1.Range("F21").Value = 4301
$ws1.Range("F25").Value = 20

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 215
$ws4.Range("F4").Value = 11
$ws4.Range("F5").Value = 6606
$ws4.Range("F9").Value = 5989
$ws4.Range("F16").Value = 92
$ws4.Range("F21").Value = 4301
$ws4.Range("F26").Value = 20
